$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.35"
$ws.Range("E2").Value = "'4.41%"
$ws.Range("D3").Value = "'35.91"
$ws.Range("E3").Value = "'15.69%"
$ws.Range("D4").Value = "'5.093"
$ws.Range("E4").Value = "'2.90%"
$ws.Range("D5").Value = "'0.07863"
$ws.Range("E5").Value = "'5.60%"
$ws.Range("D6").Value = "'2.287"
$ws.Range("E6").Value = "'2.40%"
$ws.Range("D7").Value = "'8.088"
$ws.Range("E7").Value = "'4.79%"
$ws.Range("D8").Value = "'4.006"
$ws.Range("E8").Value = "'6.92%"
$ws.Range("D9").Value = "'0.9252"
$ws.Range("E9").Value = "'0.75%"
$ws.Range("D10").Value = "'0.1005"
$ws.Range("E10").Value = "'6.82%"
$ws.Range("D11").Value = "'0.1826"
$ws.Range("E11").Value = "'6.09%"
$ws.Range("D12").Value = "'0.08689"
$ws.Range("E12").Value = "'4.43%"
$ws.Range("D13").Value = "'0.03412"
$ws.Range("E13").Value = "'5.66%"
$ws.Range("D14").Value = "'0.09904"
$ws.Range("E14").Value = "'-0.54%"
$ws.Range("D15").Value = "'0.001482"
$ws.Range("E15").Value = "'-0.65%"
$ws.Range("D16").Value = "'0.04663"
$ws.Range("E16").Value = "'2.95%"
$ws.Range("E17").Value = "'-2.04%"
$ws.Range("D18").Value = "'3.495"
$ws.Range("E18").Value = "'0.56%"
$ws.Range("E19").Value = "'-1.56%"
$ws.Range("E20").Value = "'3.16%"
$ws.Range("E21").Value = "'1.42%"
$ws.Range("D22").Value = "'4.556"
$ws.Range("E22").Value = "'9.28%"
$ws.Range("D23").Value = "'0.2233"
$ws.Range("E23").Value = "'5.56%"
$ws.Range("E24").Value = "'2.12%"
$ws.Range("D25").Value = "'0.004492"
$ws.Range("E25").Value = "'5.44%"
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("E26").Value = "'0.25%"
$ws.Range("D27").Value = "'0.0002997"
$ws.Range("E27").Value = "'-11.37%"
$ws.Range("D39").Value = "'0.01755"
$ws.Range("E39").Value = "'9.35%"
$ws.Range("D40").Value = "'0.04698"
$ws.Range("E40").Value = "'2.96%"
$ws.Range("D41").Value = "'0.007866"
$ws.Range("E41").Value = "'5.91%"
$ws.Range("D42").Value = "'0.1416"
$ws.Range("E42").Value = "'4.61%"
$ws.Range("D43").Value = "'0.008790"
$ws.Range("E43").Value = "'-10.43%"
$ws.Range("D44").Value = "'0.002215"
$ws.Range("E44").Value = "'2.94%"
$ws.Range("D45").Value = "'0.009175"
$ws.Range("E45").Value = "'-4.78%"
$ws.Range("D46").Value = "'0.00006013"
$ws.Range("E46").Value = "'-1.14%"
$ws.Range("E47").Value = "'0.25%"
$ws.Range("D48").Value = "'5.803"
$ws.Range("E48").Value = "'120.53%"
$ws.Range("D49").Value = "'0.002687"
$ws.Range("E49").Value = "'34.87%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'0.25%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'0.25%"
